# Update crypto price/volume figures for the latest data refresh.
# Values are written as literal text (matching the sheet's existing text-cell
# formatting); a leading apostrophe forces Excel to keep purely numeric-looking
# "Price" strings (single decimal point) as text instead of auto-converting them
# to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.342.33"
$ws.Range("E2").Value = "  +2.62%  "
# Row 3
$ws.Range("D3").Value = "2.308.70"
$ws.Range("E3").Value = "  +1.66%  "
# Row 4
$ws.Range("E4").Value = "  -0.01%  "
# Row 5
$ws.Range("D5").Value = "'310.88"
$ws.Range("E5").Value = "  +0.63%  "
# Row 6
$ws.Range("D6").Value = "'103.55"
$ws.Range("E6").Value = "  +6.65%  "
# Row 7
$ws.Range("E7").Value = "  +1.07%  "
# Row 8
$ws.Range("E8").Value = "  -0.03%  "
# Row 9
$ws.Range("D9").Value = "'0.530"
$ws.Range("E9").Value = "  +8.41%  "
# Row 10
$ws.Range("D10").Value = "'36.59"
$ws.Range("E10").Value = "  +4.23%  "
# Row 11
$ws.Range("D11").Value = "'52.59"
$ws.Range("E11").Value = "  +0.78%  "
# Row 12
$ws.Range("D12").Value = "'0.0812"
$ws.Range("E12").Value = "  +0.90%  "
# Row 13
$ws.Range("E13").Value = "  -1.12%  "
# Row 14
$ws.Range("D14").Value = "'6.99"
$ws.Range("E14").Value = "  +2.50%  "
# Row 15
$ws.Range("D15").Value = "2.665.70"
$ws.Range("E15").Value = "  +1.60%  "
# Row 16
$ws.Range("D16").Value = "'15.07"
$ws.Range("E16").Value = "  +2.66%  "
# Row 17
$ws.Range("D17").Value = "2.309.69"
$ws.Range("E17").Value = "  +1.87%  "
# Row 18
$ws.Range("E18").Value = "  +2.51%  "
# Row 19
$ws.Range("D19").Value = "43.248.04"
$ws.Range("E19").Value = "  +2.73%  "
# Row 20
$ws.Range("E20").Value = "  -0.27%  "
# Row 21
$ws.Range("E21").Value = "  +2.37%  "
# Row 22
$ws.Range("D22").Value = "'6.17"
$ws.Range("E22").Value = "  +3.35%  "
# Row 23
$ws.Range("D23").Value = "'68.09"
$ws.Range("E23").Value = "  +0.71%  "
# Row 24
$ws.Range("D24").Value = "'242.44"
$ws.Range("E24").Value = "  +2.61%  "
# Row 25
$ws.Range("E25").Value = "  +2.49%  "
# Row 26
$ws.Range("E26").Value = "  +0.93%  "
# Row 27
$ws.Range("E27").Value = "  +0.21%  "
# Row 28
$ws.Range("D28").Value = "'24.94"
$ws.Range("E28").Value = "  +5.99%  "
# Row 29
$ws.Range("E29").Value = "  +8.05%  "
# Row 30
$ws.Range("D30").Value = "'37.02"
$ws.Range("E30").Value = "  -0.32%  "
# Row 31
$ws.Range("D31").Value = "'9.65"
$ws.Range("E31").Value = "  +1.11%  "
# Row 32
$ws.Range("D32").Value = "'167.44"
$ws.Range("E32").Value = "  +2.27%  "
# Row 33
$ws.Range("E33").Value = "  +0.49%  "
# Row 34
$ws.Range("E34").Value = "  -0.02%  "
# Row 35
$ws.Range("D35").Value = "'18.37"
$ws.Range("E35").Value = "  +4.39%  "
# Row 36
$ws.Range("E36").Value = "  +6.77%  "
# Row 37
$ws.Range("D37").Value = "'0.0743"
$ws.Range("E37").Value = "  +1.21%  "
# Row 38
$ws.Range("E38").Value = "  -1.39%  "
# Row 39
$ws.Range("D39").Value = "'1.88"
$ws.Range("E39").Value = "  +3.35%  "
# Row 40
$ws.Range("E40").Value = "  +1.71%  "
# Row 41
$ws.Range("E41").Value = "  +6.87%  "
# Row 42
$ws.Range("E42").Value = "  +0.80%  "
# Row 43
$ws.Range("D43").Value = "'2.66"
$ws.Range("E43").Value = "  +16.93%  "
# Row 44
$ws.Range("E44").Value = "  +3.89%  "
# Row 45
$ws.Range("D45").Value = "1.987.22"
$ws.Range("E45").Value = "  +1.95%  "
# Row 46
$ws.Range("D46").Value = "'18.90"
$ws.Range("E46").Value = "  +0.54%  "
# Row 47
$ws.Range("D47").Value = "'3.05"
$ws.Range("E47").Value = "  +3.08%  "
# Row 48
$ws.Range("D48").Value = "'10.01"
$ws.Range("E48").Value = "  +2.47%  "
# Row 49
$ws.Range("D49").Value = "'55.78"
$ws.Range("E49").Value = "  +3.48%  "
# Row 50
$ws.Range("D50").Value = "'2.94"
$ws.Range("E50").Value = "  +0.81%  "
# Row 51
$ws.Range("E51").Value = "  +8.43%  "
